$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date column (A) to the new log date for each activity row
$ws.Range("A2").Value2 = 44177
$ws.Range("A3").Value2 = 44177
$ws.Range("A4").Value2 = 44177
$ws.Range("A5").Value2 = 44177

# Update the "Time Spent" column (C) with the newly recorded durations
$ws.Range("C2").Value = "00:40:11"
$ws.Range("C3").Value = "00:26:56"
$ws.Range("C4").Value = "00:01:56"
$ws.Range("C5").Value = "06:03:31"

# Move the active selection to E6
$ws.Range("E6").Select()
